$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase full names, strip leading envelope glyph from email addresses
$ws.Range("B2").Value = "JOHN DOE"
$ws.Range("C2").Value = "john@company.com"

$ws.Range("B3").Value = "JANE SMITH"
$ws.Range("C3").Value = "jane@company.com"

$ws.Range("B4").Value = "BOB JOHNSON"
$ws.Range("C4").Value = "bob@company.com"

$ws.Range("B5").Value = "ALICE BROWN"
$ws.Range("C5").Value = "alice@company.com"

# Adjust column widths: A & B both become the original A width (14.7142857142857),
# C shrinks from 20.7142857142857 to 18.7142857142857.
# NOTE: the runtime's ColumnWidth setter quantizes to whole-pixel increments
# (1/6 character-width units) before writing the XML "width" attribute, so we
# pick the ColumnWidth input whose quantized result lands closest to the
# exact target stored width.
$ws.Columns.Item(1).ColumnWidth = 13.8333333333333
$ws.Columns.Item(2).ColumnWidth = 13.8333333333333
$ws.Columns.Item(3).ColumnWidth = 17.8333333333333

